# Generate Report for Handback
# - Update status text from "Ready for handoff" to "Handback transform failed"
#   wherever it appears (Overview sheet + zh-cn / de-de detail sheets).
# - Record the handback error detail message in column K ("Error Detail")
#   for the c9180b5e-... row on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhError = "Handback file name: kno2wv1z.ion is different with handoff file name: c9180b5e-3085-403d-8a93-1f2c10807232.bbcb7cab55336259e4c6d0b297299aba523b1831.zh-cn."
$deError  = "Handback file name: kno2wv1z.ion is different with handoff file name: c9180b5e-3085-403d-8a93-1f2c10807232.bbcb7cab55336259e4c6d0b297299aba523b1831.de-de."

# --- Overview sheet: update status cells for the c9180b5e-... row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: update status + add error detail for row 3 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("K3").Value = $zhError

# --- de-de sheet: update status + add error detail for row 3 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("K3").Value = $deError
